$d = $word.ActiveDocument

# --- 1) Merge the <id>p097r_1</id> runs into a single run ---------------
$null = $d.Content.Find.Execute(
    "<id>p097r_1</id>", $true, $false, $false, $false, $false, $true,
    1, $false, "<id>p097r_1</id>", 2)

# --- 2) "ou" -> "où" : " du livre ou tu le peulx voyr" becomes
#        " du livre où tu le peulx voyr" (only the "u" becomes "u-grave")
$rng = $d.Content
$found = $rng.Find.Execute(
    "du livre ou tu le peulx voyr", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if ($found) {
    $uStart = $rng.Start + 10
    $uEnd = $uStart + 1
    $uRange = $d.Range($uStart, $uEnd)
    $uRange.Text = [char]0x00F9
    $uRange2 = $d.Range($uStart, $uStart + 1)
    $uRange2.Font.Color = -16777216
}

# --- 3) Merge the <id>p097r_2</id> runs into a single run ---------------
$null = $d.Content.Find.Execute(
    "<id>p097r_2</id>", $true, $false, $false, $false, $false, $true,
    1, $false, "<id>p097r_2</id>", 2)

# --- 4) Add a comma: "un petit le bleu" -> "un petit, le bleu" ----------
$null = $d.Content.Find.Execute(
    "sinon un petit le bleu", $true, $false, $false, $false, $false, $true,
    1, $false, "sinon un petit, le bleu", 2)

# --- 5) Merge the <id>p097r_3</id> runs into a single run ---------------
$null = $d.Content.Find.Execute(
    "<id>p097r_3</id>", $true, $false, $false, $false, $false, $true,
    1, $false, "<id>p097r_3</id>", 2)
